$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1926896.2
$ws.Range("J17").Value = 1964665.5
$ws.Range("L17").Value = 5893996.5
$ws.Range("N17").Value = -5894332.5
$ws.Range("H33").Value = 131.35715
$ws.Range("I33").Value = 133.76923
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 133.76923
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = 95.23077000000001
$ws.Range("N33").Value = -558
$ws.Range("H113").Value = 3308.7778
$ws.Range("I113").Value = 2336.25
$ws.Range("J113").Value = 4086.8
$ws.Range("K113").Value = 2336.25
$ws.Range("L113").Value = 4086.8
$ws.Range("M113").Value = 917.75
$ws.Range("N113").Value = -10594.8
$ws.Range("H132").Value = 6565
$ws.Range("I132").Value = 6758.636
$ws.Range("K132").Value = 20275.908
$ws.Range("M132").Value = -17745.908
$ws.Range("H138").Value = 1801.277
$ws.Range("J138").Value = 2402.7273
$ws.Range("L138").Value = 7208.1819
$ws.Range("N138").Value = -17488.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2508.361
$ws.Range("I61").Value = 2163.0952
$ws.Range("J61").Value = 2991.7334
$ws.Range("K61").Value = 2163.0952
$ws.Range("L61").Value = 2991.7334
$ws.Range("M61").Value = -1951.0952
$ws.Range("N61").Value = -3415.7334
$ws.Range("H74").Value = 71432300
$ws.Range("I74").Value = 90913440
$ws.Range("J74").Value = 1471.3334
$ws.Range("K74").Value = 90913440
$ws.Range("L74").Value = 1471.3334
$ws.Range("M74").Value = -90912566
$ws.Range("N74").Value = -3219.3334
$ws.Range("H77").Value = 71432300
$ws.Range("I77").Value = 90913440
$ws.Range("J77").Value = 1471.3334
$ws.Range("K77").Value = 454567200
$ws.Range("L77").Value = 7356.666999999999
$ws.Range("M77").Value = -454562832
$ws.Range("N77").Value = -16092.667
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H96").Value = 18366.375
$ws.Range("J96").Value = 18366.375
$ws.Range("L96").Value = 18366.375
$ws.Range("N96").Value = -23858.375
$ws.Range("H97").Value = 923.625
$ws.Range("I97").Value = 1248.091
$ws.Range("J97").Value = 209.8
$ws.Range("K97").Value = 1248.091
$ws.Range("L97").Value = 209.8
$ws.Range("M97").Value = -752.0909999999999
$ws.Range("N97").Value = -1201.8
$ws.Range("H122").Value = 2018.9
$ws.Range("I122").Value = 2385.8462
$ws.Range("J122").Value = 1337.4286
$ws.Range("K122").Value = 7157.5386
$ws.Range("L122").Value = 4012.2858
$ws.Range("M122").Value = -4707.5386
$ws.Range("N122").Value = -8912.2858
$ws.Range("H132").Value = 15357.216
$ws.Range("I132").Value = 1655.64
$ws.Range("J132").Value = 43902.168
$ws.Range("K132").Value = 4966.92
$ws.Range("L132").Value = 131706.504
$ws.Range("M132").Value = -2436.92
$ws.Range("N132").Value = -136766.504
$ws.Range("H136").Value = 2508.361
$ws.Range("I136").Value = 2163.0952
$ws.Range("J136").Value = 2991.7334
$ws.Range("K136").Value = 6489.285600000001
$ws.Range("L136").Value = 8975.200199999999
$ws.Range("M136").Value = -3939.285600000001
$ws.Range("N136").Value = -14075.2002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1079.5294
$ws.Range("I94").Value = 932.09753
$ws.Range("K94").Value = 932.09753
$ws.Range("M94").Value = -481.09753
$ws.Range("H105").Value = 4548100
$ws.Range("I105").Value = 2180
$ws.Range("J105").Value = 8336366.5
$ws.Range("K105").Value = 2180
$ws.Range("L105").Value = 8336366.5
$ws.Range("M105").Value = -433
$ws.Range("N105").Value = -8339860.5
$ws.Range("H107").Value = 1295.5834
$ws.Range("I107").Value = 754.25
$ws.Range("J107").Value = 2378.25
$ws.Range("K107").Value = 754.25
$ws.Range("L107").Value = 2378.25
$ws.Range("M107").Value = 1165.75
$ws.Range("N107").Value = -6218.25
$ws.Range("H134").Value = 35210.53
$ws.Range("I134").Value = 44629.48
$ws.Range("K134").Value = 133888.44
$ws.Range("M134").Value = -131353.44

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19882.834
$ws.Range("I31").Value = 26222.385
$ws.Range("J31").Value = 3400
$ws.Range("K31").Value = 26222.385
$ws.Range("L31").Value = 3400
$ws.Range("M31").Value = -25927.385
$ws.Range("N31").Value = -3990
$ws.Range("H34").Value = 19882.834
$ws.Range("I34").Value = 26222.385
$ws.Range("J34").Value = 3400
$ws.Range("K34").Value = 26222.385
$ws.Range("L34").Value = 3400
$ws.Range("M34").Value = -26020.385
$ws.Range("N34").Value = -3804
$ws.Range("H59").Value = 23736.842
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").Value = ""
$ws.Range("H94").Value = 9937.4
$ws.Range("J94").Value = 12910.571
$ws.Range("L94").Value = 12910.571
$ws.Range("N94").Value = -13812.571
$ws.Range("H99").Value = 11908244
$ws.Range("I99").Value = 2876.1155
$ws.Range("J99").Value = 31254468
$ws.Range("K99").Value = 2876.1155
$ws.Range("L99").Value = 31254468
$ws.Range("M99").Value = -1378.1155
$ws.Range("N99").Value = -31257464
$ws.Range("H105").Value = 7813698
$ws.Range("J105").Value = 1721.5
$ws.Range("L105").Value = 1721.5
$ws.Range("N105").Value = -5215.5
$ws.Range("H122").Value = 2493.6667
$ws.Range("J122").Value = 1457
$ws.Range("L122").Value = 4371
$ws.Range("N122").Value = -9271
$ws.Range("H126").Value = 11908244
$ws.Range("I126").Value = 2876.1155
$ws.Range("J126").Value = 31254468
$ws.Range("K126").Value = 8628.3465
$ws.Range("L126").Value = 93763404
$ws.Range("M126").Value = -6158.3465
$ws.Range("N126").Value = -93768344

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2307786.5
$ws.Range("I4").Value = 92.3
$ws.Range("J4").Value = 10000100
$ws.Range("K4").Value = 276.9
$ws.Range("L4").Value = 30000300
$ws.Range("M4").Value = -164.9
$ws.Range("N4").Value = -30000524
$ws.Range("H131").Value = 766.54
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 766.54
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2299.62
$ws.Range("M131").Value = ""
$ws.Range("N131").Value = -12379.62
$ws.Range("H139").Value = 1755.2
$ws.Range("I139").Value = 1148.2667
$ws.Range("J139").Value = 2665.6
$ws.Range("K139").Value = 3444.800099999999
$ws.Range("L139").Value = 7996.799999999999
$ws.Range("M139").Value = 1695.199900000001
$ws.Range("N139").Value = -18276.8
$ws.Range("H140").Value = 1673.5714
$ws.Range("I140").Value = 1404.4445
$ws.Range("K140").Value = 4213.333500000001
$ws.Range("M140").Value = 966.6664999999994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2295
$ws.Range("I113").Value = 1792.9286
$ws.Range("K113").Value = 1792.9286
$ws.Range("M113").Value = 377.0714
$ws.Range("H126").Value = 3971.3713
$ws.Range("I126").Value = 3049.9092
$ws.Range("K126").Value = 9149.7276
$ws.Range("M126").Value = -6679.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3048.2
$ws.Range("I61").Value = 1292.25
$ws.Range("K61").Value = 1292.25
$ws.Range("M61").Value = -1090.25
$ws.Range("H93").Value = 2256.9333
$ws.Range("I93").Value = 1988.1538
$ws.Range("J93").Value = 4004
$ws.Range("K93").Value = 1988.1538
$ws.Range("L93").Value = 4004
$ws.Range("M93").Value = -740.1538
$ws.Range("N93").Value = -6500
$ws.Range("H100").Value = 2231.92
$ws.Range("I100").Value = 1883.1111
$ws.Range("J100").Value = 3128.8572
$ws.Range("K100").Value = 1883.1111
$ws.Range("L100").Value = 3128.8572
$ws.Range("M100").Value = -1342.1111
$ws.Range("N100").Value = -4210.8572
$ws.Range("H113").Value = 3048.2
$ws.Range("I113").Value = 1292.25
$ws.Range("K113").Value = 1292.25
$ws.Range("M113").Value = 877.75
$ws.Range("H122").Value = 1404139.4
$ws.Range("J122").Value = 5300.8
$ws.Range("L122").Value = 15902.4
$ws.Range("N122").Value = -20802.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1269.56
$ws.Range("I126").Value = 1408.6
$ws.Range("K126").Value = 4225.799999999999
$ws.Range("M126").Value = -1755.799999999999
